$wb = $excel.ActiveWorkbook

# Update the "OFF" sheet (Wild Card round stats added to row 2: "H")
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 499
$wsOff.Range("C2").Value = 350
$wsOff.Range("D2").Value = 106
$wsOff.Range("E2").Value = 42
$wsOff.Range("G2").Value = 5

# Update the "DEF" sheet (Wild Card round stats added to row 2: "H")
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 432
$wsDef.Range("C2").Value = 319
$wsDef.Range("D2").Value = 103
$wsDef.Range("E2").Value = 55
